$wb = $excel.ActiveWorkbook

# --- Sheet: Recommandations ---
$ws = $wb.Worksheets.Item("Recommandations")

$ws.Cells.Item(2,4).Value = 83456.25999999999
$ws.Cells.Item(3,4).Value = 78075
$ws.Cells.Item(4,4).Value = 72655
$ws.Cells.Item(5,4).Value = 65354.52
$ws.Cells.Item(6,4).Value = 61500
$ws.Cells.Item(7,4).Value = 60440
$ws.Cells.Item(8,4).Value = 58900
$ws.Cells.Item(9,4).Value = 54695
$ws.Cells.Item(10,4).Value = 48985
$ws.Cells.Item(11,4).Value = 42550
$ws.Cells.Item(12,4).Value = 36587.3
$ws.Cells.Item(13,4).Value = 32595.13
$ws.Cells.Item(14,4).Value = 22926.22
$ws.Cells.Item(16,4).Value = 15077.08
$ws.Cells.Item(17,4).Value = 13787.5
$ws.Cells.Item(18,4).Value = 11873.69
$ws.Cells.Item(19,4).Value = 11869.88
$ws.Cells.Item(20,4).Value = 11333.11
$ws.Cells.Item(21,4).Value = 11186.45
$ws.Cells.Item(22,4).Value = 10927.83
$ws.Cells.Item(23,4).Value = 10739.59
$ws.Cells.Item(24,4).Value = 9663.969999999999
$ws.Cells.Item(25,4).Value = 9473.24
$ws.Cells.Item(30,4).Value = 78.3
$ws.Cells.Item(30,7).Value = '✅ Renforcer'
$ws.Cells.Item(31,1).Value = 'TOTALENERGIES MARKETING CI (TTLC)'
$ws.Cells.Item(31,2).Value = 24
$ws.Cells.Item(31,3).Value = 19
$ws.Cells.Item(31,4).Value = 50.97
$ws.Cells.Item(31,5).Value = 4.45
$ws.Cells.Item(31,7).Value = '✅ Renforcer'
$ws.Cells.Item(32,1).Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws.Cells.Item(32,2).Value = 12
$ws.Cells.Item(32,3).Value = 2
$ws.Cells.Item(32,4).Value = 48.5
$ws.Cells.Item(32,5).Value = -3.77
$ws.Cells.Item(32,7).Value = '➖ Neutre'
$ws.Cells.Item(33,1).Value = 'BANK OF AFRICA SENEGAL (BOAS)'
$ws.Cells.Item(33,3).Value = 5
$ws.Cells.Item(33,4).Value = 47.13
$ws.Cells.Item(33,5).Value = 3.92
$ws.Cells.Item(33,7).Value = 'Non évalué'
$ws.Cells.Item(34,1).Value = 'SAPH CI (SPHC)'
$ws.Cells.Item(34,2).Value = 15
$ws.Cells.Item(34,3).Value = 9
$ws.Cells.Item(34,4).Value = 46.42
$ws.Cells.Item(34,5).Value = 3.33
$ws.Cells.Item(34,7).Value = '➖ Neutre'
$ws.Cells.Item(36,2).Value = 8
$ws.Cells.Item(36,4).Value = 42.51
$ws.Cells.Item(38,1).Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws.Cells.Item(38,2).Value = 8
$ws.Cells.Item(38,4).Value = 27.56
$ws.Cells.Item(38,5).Value = 5.82
$ws.Cells.Item(38,7).Value = '✅ Renforcer'
$ws.Cells.Item(39,1).Value = 'BANK OF AFRICA CI (BOAC)'
$ws.Cells.Item(39,2).Value = 7
$ws.Cells.Item(39,4).Value = 26.08
$ws.Cells.Item(39,5).Value = -1.88
$ws.Cells.Item(39,7).Value = '👀 À surveiller'
$ws.Cells.Item(40,1).Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws.Cells.Item(40,2).Value = 24
$ws.Cells.Item(40,3).Value = 22
$ws.Cells.Item(40,4).Value = 24.09
$ws.Cells.Item(40,5).Value = 6.67
$ws.Cells.Item(40,7).Value = '✅ Renforcer'
$ws.Cells.Item(41,1).Value = 'BANK OF AFRICA BN (BOAB)'
$ws.Cells.Item(41,2).Value = 6
$ws.Cells.Item(41,3).Value = 3
$ws.Cells.Item(41,4).Value = 23.86
$ws.Cells.Item(41,5).Value = -1.9
$ws.Cells.Item(41,7).Value = '➖ Neutre'
$ws.Cells.Item(42,1).Value = 'SMB CI (SMBC)'
$ws.Cells.Item(42,2).Value = 9
$ws.Cells.Item(42,3).Value = 8
$ws.Cells.Item(42,4).Value = 23.09
$ws.Cells.Item(42,5).Value = -2.07
$ws.Cells.Item(42,7).Value = 'Non évalué'
$ws.Cells.Item(43,1).Value = 'CFAO MOTORS CI (CFAC)'
$ws.Cells.Item(43,2).Value = 11
$ws.Cells.Item(43,4).Value = 21.59
$ws.Cells.Item(43,5).Value = 6.98
$ws.Cells.Item(43,7).Value = '➖ Neutre'
$ws.Cells.Item(44,1).Value = 'BANK OF AFRICA ML (BOAM)'
$ws.Cells.Item(44,2).Value = 13
$ws.Cells.Item(44,3).Value = 9
$ws.Cells.Item(44,4).Value = 18.35
$ws.Cells.Item(44,5).Value = 3.2
$ws.Cells.Item(44,7).Value = 'Non évalué'
$ws.Cells.Item(45,3).Value = 9
$ws.Cells.Item(45,4).Value = 18.31
$ws.Cells.Item(46,1).Value = 'ORAGROUP TOGO (ORGT)'
$ws.Cells.Item(46,2).Value = 7
$ws.Cells.Item(46,3).Value = 4
$ws.Cells.Item(46,4).Value = 17.5
$ws.Cells.Item(46,5).Value = 5.63
$ws.Cells.Item(46,7).Value = '➖ Neutre'
$ws.Cells.Item(57,4).Value = -2.03
$ws.Cells.Item(58,3).Value = 9
$ws.Cells.Item(58,4).Value = -8.76
$ws.Cells.Item(64,1).Value = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$ws.Cells.Item(64,2).Value = 1
$ws.Cells.Item(64,3).Value = 11
$ws.Cells.Item(64,4).Value = -16.75
$ws.Cells.Item(64,5).Value = -1.14
$ws.Cells.Item(64,7).Value = 'Non évalué'
$ws.Cells.Item(65,1).Value = 'SICABLE CI (CABC)'
$ws.Cells.Item(65,2).Value = 17
$ws.Cells.Item(65,3).Value = 28
$ws.Cells.Item(65,4).Value = -19.38
$ws.Cells.Item(65,5).Value = -4.41
$ws.Cells.Item(65,7).Value = '⚠️ Risque de décrochage'
$ws.Cells.Item(66,1).Value = 'SETAO CI (STAC)'
$ws.Cells.Item(66,2).Value = 22
$ws.Cells.Item(66,3).Value = 25
$ws.Cells.Item(66,4).Value = -19.63
$ws.Cells.Item(66,5).Value = -2.5
$ws.Cells.Item(66,7).Value = '👀 À surveiller'
$ws.Cells.Item(67,1).Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws.Cells.Item(67,2).Value = 7
$ws.Cells.Item(67,3).Value = 14
$ws.Cells.Item(67,4).Value = -23.32
$ws.Cells.Item(67,5).Value = 4.92
$ws.Cells.Item(67,7).Value = '✅ Renforcer'
$ws.Cells.Item(68,1).Value = 'AIR LIQUIDE CI (SIVC)'
$ws.Cells.Item(68,2).Value = 3
$ws.Cells.Item(68,3).Value = 12
$ws.Cells.Item(68,4).Value = -27.61
$ws.Cells.Item(68,5).Value = -5.45
$ws.Cells.Item(68,7).Value = '➖ Neutre'
$ws.Cells.Item(72,2).Value = 7
$ws.Cells.Item(72,4).Value = -43.72

# --- Sheet: Top_YTD ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Cells.Item(5,2).Value = 105.43
$ws2.Cells.Item(6,1).Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws2.Cells.Item(6,2).Value = 59.76
$ws2.Cells.Item(7,1).Value = 'BANK OF AFRICA SENEGAL (BOAS)'
$ws2.Cells.Item(7,2).Value = 57.81
$ws2.Cells.Item(8,1).Value = 'TOTALENERGIES MARKETING CI (TTLC)'
$ws2.Cells.Item(8,2).Value = 54.89
$ws2.Cells.Item(9,1).Value = 'SAPH CI (SPHC)'
$ws2.Cells.Item(9,2).Value = 54.75
$ws2.Cells.Item(10,1).Value = 'BICI CI (BICC)'
$ws2.Cells.Item(10,2).Value = 49.27
$ws2.Cells.Item(11,1).Value = 'UNIWAX CI (UNXC)'
$ws2.Cells.Item(11,2).Value = 45.79
